$d = $word.ActiveDocument

# 1. Agreement fix: "... pour la base de données exporté, avec ..." ->
#    "... pour la base de données exportée, avec ..." (insert a trailing "e").
$ok1 = $d.Content.Find.Execute(
    "exporté, avec l’instruction CREATE DATABASE",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "exportée, avec l’instruction CREATE DATABASE",
    2)
"Typo fix applied: $ok1"

# 2. Drop the trailing bordered paragraph entirely (text + paragraph mark + its
#    own pPr/border), so the previous paragraph now runs straight into the
#    section break.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("noter qu") -and $p.Range.Text.Contains("aucun administrateur")) {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
    "Trailing paragraph removed"
}
else {
    "Trailing paragraph NOT found"
}
